$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates: reference => new value
$updates = [ordered]@{
    'D2' = '60.856.24'
    'E2' = '  +2.67%  '
    'D3' = '2.610.32'
    'E3' = '  +1.37%  '
    'E4' = '  +0.03%  '
    'D5' = '579.66'
    'E5' = '  +4.30%  '
    'D6' = '144.01'
    'E6' = '  +1.62%  '
    'E7' = '  -0.25%  '
    'D8' = '0.602'
    'E8' = '  +0.64%  '
    'D9' = '2.637.44'
    'E9' = '  +2.18%  '
    'E10' = '  -2.95%  '
    'E11' = '  +2.44%  '
    'D12' = '0.158'
    'E12' = '  -3.70%  '
    'D13' = '0.372'
    'E13' = '  +5.99%  '
    'D14' = '3.079.90'
    'E14' = '  +1.57%  '
    'D15' = '60.846.98'
    'E15' = '  +2.66%  '
    'E16' = '  +2.02%  '
    'E17' = '  +4.46%  '
    'D18' = '2.624.62'
    'E18' = '  +1.81%  '
    'E19' = '  +9.41%  '
    'E20' = '  +3.02%  '
    'D21' = '350.75'
    'E21' = '  +3.81%  '
    'D22' = '6.95'
    'E22' = '  +7.53%  '
    'E23' = '  +0.04%  '
    'E24' = '  +7.82%  '
    'D25' = '63.29'
    'E25' = '  +1.10%  '
    'E26' = '  -0.19%  '
    'E27' = '  +0.61%  '
    'E28' = '  +7.25%  '
    'E29' = '  +3.65%  '
    'E30' = '  +8.93%  '
    'D31' = '6.36'
    'E31' = '  +2.77%  '
    'E32' = '  -0.13%  '
    'D33' = '163.08'
    'E33' = '  +2.35%  '
    'E34' = '  +2.78%  '
    'B35' = 'Fetch.AI'
    'C35' = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    'D35' = '1.01'
    'E35' = '  +13.40%  '
    'B36' = 'NEARProtocol'
    'C36' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D36' = '4.30'
    'E36' = '  +5.43%  '
    'E37' = '  +6.22%  '
    'E38' = '  +10.08%  '
    'E39' = '  +1.44%  '
    'D40' = '3.91'
    'E40' = '  +6.49%  '
    'D41' = '308.37'
    'E41' = '  +6.68%  '
    'D42' = '0.849'
    'E42' = '  -0.21%  '
    'D43' = '134.17'
    'E43' = '  -3.10%  '
    'D44' = '20.43'
    'E44' = '  +9.61%  '
    'E45' = '  -0.40%  '
    'B46' = 'RenderToken'
    'C46' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D46' = '5.05'
    'E46' = '  +12.00%  '
    'E47' = '  +2.94%  '
    'B48' = 'EnergySwap'
    'C48' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D48' = '19.92'
    'E48' = '  +5.57%  '
    'E49' = '  +1.19%  '
    'D50' = '0.0551'
    'E50' = '  +4.15%  '
    'E51' = '  +3.95%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    if ($ref[0] -eq "D") {
        # Price column: force text so values like "579.66" are not
        # auto-converted to numbers (matches original inlineStr formatting).
        $cell.NumberFormat = "@"
        $cell.Value = $updates[$ref]
        $cell.Style = "Normal"
    } else {
        $cell.Value = $updates[$ref]
    }
}
